$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Upload the following to the " -> "Upload the following 6 files to the "
# ---------------------------------------------------------------------------
$p42 = $d.Paragraphs.Item(42).Range
$ok42 = $p42.Find.Execute(
    "Upload the following to the ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Upload the following 6 files to the ", 2)

# ---------------------------------------------------------------------------
# 2) "A zip file containing the four files" -> "The four files"
# ---------------------------------------------------------------------------
$p44 = $d.Paragraphs.Item(44).Range
$ok44 = $p44.Find.Execute(
    "A zip file containing the four files", $true, $false, $false, $false,
    $false, $true, 1, $false, "The four files", 2)

# ---------------------------------------------------------------------------
# 3) " code review of your own code." ->
#    " code review of your code with the “Prod” column filled in by you."
#    (only "“Prod” column filled in by you" stays underlined; the trailing
#    period is not underlined)
# ---------------------------------------------------------------------------
$p45a = $d.Paragraphs.Item(45).Range
$ok45a = $p45a.Find.Execute(
    " code review ", $true, $false, $false, $false, $false, $true, 1,
    $false, " code review of your code with the ", 2)

$p45b = $d.Paragraphs.Item(45).Range
$ok45b = $p45b.Find.Execute(
    "of your own code.", $true, $false, $false, $false, $false, $true, 1,
    $false, "“Prod” column filled in by you.", 2)

# Un-underline just the trailing period (it was carried over as underlined
# from the run it replaced).
$p45c = $d.Paragraphs.Item(45).Range
$lastPeriod = $d.Range($p45c.End - 2, $p45c.End - 1)
$lastPeriod.Font.Underline = 0

# ---------------------------------------------------------------------------
# 4) Move the "_GoBack" bookmark: it currently wraps (empty) in paragraph 44;
#    it should instead start right before paragraph 42's text and end right
#    after paragraph 45's text.
# ---------------------------------------------------------------------------
$spanStart = $d.Paragraphs.Item(42).Range.Start
$spanEnd = $d.Paragraphs.Item(45).Range.End
$span = $d.Range($spanStart, $spanEnd)
$d.Bookmarks.Add("_GoBack", $span) | Out-Null

Write-Output "done: ok42=$ok42 ok44=$ok44 ok45a=$ok45a ok45b=$ok45b"
